$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(2).HorizontalAlignment = -4131
Write-Host "done"
